$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "lelijke_excel_dates"

# New "exceldate" numeric values in column E (rows 2-4)
$ws.Range("E2").Value = 41080
$ws.Range("E4").Value = 41081

# E3 has no value but still carries the column's cell formatting, so touch
# its border (a visual no-op) to force the cell to materialize in the sheet.
$ws.Range("E3").Borders.LineStyle = -4142

# Move the active selection to F1 (mirrors the saved selection state in the
# edited workbook)
$ws.Range("F1").Select()
